$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'27.720.05"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.57%  '
$ws.Range('D3').Value = "'1.639.70"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.44%  '
$ws.Range('D5').Value = "'212.76"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.08%  '
$ws.Range('E6').Value = '  -2.31%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = "'23.30"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.07%  '
$ws.Range('D9').Value = "'0.261"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.52%  '
$ws.Range('E10').Value = '  +0.06%  '
$ws.Range('E11').Value = '  +0.19%  '
$ws.Range('D12').Value = "'1.871.51"
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Value = "'1.656.42"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.64%  '
$ws.Range('E14').Value = '  +0.50%  '
$ws.Range('E16').Value = '  +0.53%  '
$ws.Range('D17').Value = "'27.681.56"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.54%  '
$ws.Range('D18').Value = "'230.33"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.27%  '
$ws.Range('D19').Value = "'7.70"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.05%  '
$ws.Range('E20').Value = '  -0.10%  '
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('D22').Value = "'4.31"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('D23').Value = "'10.25"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.40%  '
$ws.Range('E24').Value = '  +3.28%  '
$ws.Range('D25').Value = "'151.04"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.36%  '
$ws.Range('E26').Value = '  -0.76%  '
$ws.Range('E27').Value = '  -0.75%  '
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('E29').Value = '  +0.12%  '
$ws.Range('E30').Value = '  +0.22%  '
$ws.Range('D31').Value = "'0.0488"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.35%  '
$ws.Range('E32').Value = '  +0.26%  '
$ws.Range('D33').Value = "'1.464.06"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.76%  '
$ws.Range('E34').Value = '  -2.01%  '
$ws.Range('E35').Value = '  -1.98%  '
$ws.Range('D37').Value = "'0.569"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('D38').Value = "'0.882"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.16%  '
$ws.Range('E39').Value = '  +0.33%  '
$ws.Range('D40').Value = "'0.894"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +9.17%  '
$ws.Range('D41').Value = "'69.04"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +6.11%  '
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('E43').Value = '  -2.03%  '
$ws.Range('E44').Value = '  +1.26%  '
$ws.Range('D45').Value = "'2.45"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.82%  '
$ws.Range('E46').Value = '  -0.73%  '
$ws.Range('D47').Value = "'1.781.61"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E48').Value = '  +3.26%  '
$ws.Range('D49').Value = "'87.16"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.09%  '
$ws.Range('E50').Value = '  -1.20%  '
$ws.Range('D51').Value = "'0.0994"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.01%  '
